$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new log entry row (row 10) following the same pattern as row 9
$ws.Range("A10").Value = "2/13/2020jaclemon"
$ws.Range("C10").Value = "45 minutes"
$ws.Range("D10").Value = "Used Clion to diagnose illegal command error"

# Keep selection consistent with the existing sheet view state
$ws.Range("D10").Select()
